{"js": "// Target change (paragraph containing \"...must be a bad drug for the\n// master, if he is red-haired \\u2014 and for the apprentice too.\"):\n//   1. Insert a new run reading \"both \" immediately before \"for \" (so the\n//      sentence reads \"...a bad drug both for the master...\").\n//   2. Collapse the three runs \"is red-haired \\u2014\" + \" \" + \"and for \"\n//      into a single run reading \"is red-haired, and for \".\n\nconst body = context.document.body;\n\n// --- Step 1: insert \"both \" as its own run -------------------------------\n// Insert the text first (Office.js/Word happily coalesces it into the\n// preceding run since the formatting is identical at that point).\nconst forResults = body.search(\"for the master\", { matchCase: true, matchWholeWord: false });\nforResults.load(\"items\");\nawait context.sync();\nif (forResults.items.length !== 1) {\n  throw new Error(`expected exactly one \"for the master\" match, found ${forResults.items.length}`);\n}\nconst forRange = forResults.items[0];\nconst insertionPoint = forRange.getRange(Word.RangeLocation.start);\ninsertionPoint.insertText(\"both \", Word.InsertLocation.before);\nawait context.sync();\n\n// Re-locate the freshly inserted text with an independent search (using the\n// direct return value of insertText for further property writes triggers\n// document-wide formatting side effects in this host, so we avoid it) and\n// toggle a character property on/off. That forces the run to remain a\n// distinct <w:r> (matching rPr/color/rtl) instead of merging back into its\n// neighbor.\nconst bothResults = body.search(\"both \", { matchCase: true, matchWholeWord: false });\nbothResults.load(\"items\");\nawait context.sync();\nif (bothResults.items.length !== 1) {\n  throw new Error(`expected exactly one \"both \" match, found ${bothResults.items.length}`);\n}\nconst bothRange = bothResults.items[0];\nbothRange.font.bold = true;\nawait context.sync();\nbothRange.font.bold = false;\nawait context.sync();\n\n// --- Step 2: merge \"is red-haired \u2014 and for \" into a single run ----------\nconst dashResults = body.search(\"is red-haired \u2014 and for \", { matchCase: true, matchWholeWord: false });\ndashResults.load(\"items\");\nawait context.sync();\nif (dashResults.items.length !== 1) {\n  throw new Error(`expected exactly one \"is red-haired \u2014 and for \" match, found ${dashResults.items.length}`);\n}\ndashResults.items[0].insertText(\"is red-haired, and for \", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Target change (paragraph containing \"...must be a bad drug for the\n# master, if he is red-haired - and for the apprentice too.\"):\n#   1. Insert a new run reading \"both \" immediately before \"for \" (so the\n#      sentence reads \"...a bad drug both for the master...\").\n#   2. Collapse the three runs \"is red-haired -\" + \" \" + \"and for \" into a\n#      single run reading \"is red-haired, and for \".\n\n$d = $word.ActiveDocument\n\n# --- Step 1: insert \"both \" as its own run --------------------------------\n$findRange = $d.Content\n$find = $findRange.Find\n$find.Text = \"for the master\"\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"could not find 'for the master'\"\n}\n\n$insertRange = $findRange.Duplicate\n$insertRange.Collapse(1)   # wdCollapseStart\n$insertRange.Text = \"both \"\n\n# Toggling a character property forces this new text to stay in its own\n# <w:r> instead of silently merging back into the identically formatted\n# run that precedes it (\"... be a bad drug \").\n$insertRange.Bold = 1\n$insertRange.Bold = 0\n\n# --- Step 2: merge \"is red-haired - and for \" into a single run ----------\n$dash = [char]8212\n$replaceRange = $d.Content\n$rfind = $replaceRange.Find\n$rfind.Text = \"is red-haired \" + $dash + \" and for \"\n$rfind.Replacement.Text = \"is red-haired, and for \"\n$replaced = $rfind.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2)\n\nif (-not $replaced) {\n    throw \"could not find/replace 'is red-haired - and for '\"\n}\n"}
